$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
Write-Host "before:" $tbl.Style
$tbl.ApplyStyle("{8E51A6B4-FEBA-4AE3-BF02-5BE07C659A7C}")
Write-Host "after:" $tbl.Style
